# New crime data collected - update weekly CompStat figures (112th Precinct)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# ---------------------------------------------------------------------------
# Header text: report volume/number and the week-covering date range
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
$ws.Range("G14").Value = 1

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 28
$ws.Range("J16").Value = 28
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 55.555555555555
$ws.Range("M16").Value = -28.205128205128
$ws.Range("N16").Value = -89.189189189189

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 200
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 71.428571428571
$ws.Range("I17").Value = 36
$ws.Range("J17").Value = 35
$ws.Range("K17").Value = 2.857142857142
$ws.Range("L17").Value = 71.428571428571
$ws.Range("M17").Value = 89.473684210526
$ws.Range("N17").Value = -25

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -22.222222222222
$ws.Range("I18").Value = 47
$ws.Range("J18").Value = 41
$ws.Range("K18").Value = 14.634146341463
$ws.Range("L18").Value = 42.424242424242
$ws.Range("M18").Value = 9.302325581395
$ws.Range("N18").Value = -90.505050505050

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -35.135135135135
$ws.Range("I19").Value = 167
$ws.Range("J19").Value = 198
$ws.Range("K19").Value = -15.656565656565
$ws.Range("L19").Value = 30.46875
$ws.Range("M19").Value = 21.897810218978
$ws.Range("N19").Value = -52.691218130311

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.  (D20/E20 switch from the "N/A" placeholders to real numbers)
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 3

$ws.Range("F20").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 5

$ws.Range("H20").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = -40

$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 40
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = 37.931034482758
$ws.Range("L20").Value = 122.222222222222
$ws.Range("M20").Value = -6.976744186046
$ws.Range("N20").Value = -96.894409937888

# ---------------------------------------------------------------------------
# Row 21 - TOTAL (Part 1 violent/property subtotal)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 5
$ws.Range("I21").Value = 319
$ws.Range("J21").Value = 338
$ws.Range("K21").Value = -5.621301775147
$ws.Range("L21").Value = 44.343891402714
$ws.Range("M21").Value = 13.120567375886
$ws.Range("N21").Value = -86.963628933387

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 11
$ws.Range("J22").Value = 13
$ws.Range("K22").Value = -15.384615384615
$ws.Range("L22").Value = 120
$ws.Range("M22").Value = 37.5

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = -52.631578947368
$ws.Range("F24").Value = 110
$ws.Range("G24").Value = 142
$ws.Range("H24").Value = -22.535211267605
$ws.Range("I24").Value = 584
$ws.Range("J24").Value = 676
$ws.Range("K24").Value = -13.609467455621
$ws.Range("L24").Value = 20.164609053497
$ws.Range("M24").Value = 64.971751412429

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = 66.666666666666
$ws.Range("I25").Value = 89
$ws.Range("J25").Value = 71
$ws.Range("K25").Value = 25.352112676056
$ws.Range("L25").Value = 41.269841269841
$ws.Range("M25").Value = 20.270270270270

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*  (D26/E26 switch from real numbers to the "N/A" placeholders)
# ---------------------------------------------------------------------------
$ws.Range("D26").Value = "'0"
$ws.Range("C26").Copy()
$ws.Range("D26").PasteSpecial(-4122)

$ws.Range("E26").Value = "***.*"
$ws.Range("C26").Copy()
$ws.Range("E26").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 11
$ws.Range("K27").Value = 22.222222222222
$ws.Range("L27").Value = -21.428571428571

# ---------------------------------------------------------------------------
# Row 28 - Shooting Vic.
# ---------------------------------------------------------------------------
$ws.Range("G28").Value = 1

# ---------------------------------------------------------------------------
# Row 29 - Shooting Inc.
# ---------------------------------------------------------------------------
$ws.Range("G29").Value = 1
